$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.179822325706482
$ws.Range("B1").Value = 2.148978471755981
$ws.Range("C1").Value = 4.370490074157715
$ws.Range("D1").Value = 2.778959274291992
$ws.Range("E1").Value = 1.221700072288513
